$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Append 4 new rows at the bottom of the sheet (rows 67-69), copying the
#    formatting used by the preceding rows (A: s=4, B: s=6 style family)
$ws.Range("A66:B66").Copy($ws.Range("A67:B69"))
$ws.Rows(67).RowHeight = 15.75
$ws.Rows(68).RowHeight = 15.75
$ws.Rows(69).RowHeight = 15.75

$ws.Range("A67").Value = "Quá trình công tác"
$ws.Range("B67").Value = "Working Progress"

$ws.Range("A68").Value = "Quá trình tập sự"
$ws.Range("B68").Value = "Probation Progress"

# 2. Update existing row 15 (Chức danh / Job title row): English text
#    "Job title, Position" -> "Job title"
$ws.Range("B15").Value = "Job title"

$ws.Range("A69").Value = "Diễn biến lương"
$ws.Range("B69").Value = "Wage Changes"

# 3. Update the view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Range("B70").Select()
